$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 14287198
$ws.Range("I43").Value = 25001248
$ws.Range("K43").Value = 25001248
$ws.Range("M43").Value = -25001179

$ws.Range("H53").Value = 249
$ws.Range("I53").Value = 258.3125
$ws.Range("K53").Value = 258.3125
$ws.Range("M53").Value = 378.6875

$ws.Range("H131").Value = 689.6667
$ws.Range("I131").Value = 689.6667
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 2069.0001
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = 2970.9999
$ws.Range("N131").Value = $null

$ws.Range("H132").Value = 18587.416
$ws.Range("I132").Value = 18587.416
$ws.Range("K132").Value = 55762.24800000001
$ws.Range("M132").Value = -53232.24800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5558.6
$ws.Range("J63").Value = 10999
$ws.Range("L63").Value = 10999
$ws.Range("N63").Value = -12371

$ws.Range("H66").Value = 5558.6
$ws.Range("J66").Value = 10999
$ws.Range("L66").Value = 54995
$ws.Range("N66").Value = -61859

$ws.Range("H86").Value = 80000
$ws.Range("J86").Value = 80000
$ws.Range("L86").Value = 80000
$ws.Range("N86").Value = -82372

$ws.Range("H88").Value = 1392.0769
$ws.Range("J88").Value = 778.6
$ws.Range("L88").Value = 778.6
$ws.Range("N88").Value = -1590.6

$ws.Range("H89").Value = 80000
$ws.Range("J89").Value = 80000
$ws.Range("L89").Value = 240000
$ws.Range("N89").Value = -251856

$ws.Range("H91").Value = 1392.0769
$ws.Range("J91").Value = 778.6
$ws.Range("L91").Value = 778.6
$ws.Range("N91").Value = -3586.6

$ws.Range("H128").Value = 120000
$ws.Range("J128").Value = 120000
$ws.Range("L128").Value = 120000
$ws.Range("N128").Value = -129960

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 23799.8
$ws.Range("J76").Value = 23799.8
$ws.Range("L76").Value = 23799.8
$ws.Range("N76").Value = -24429.8

$ws.Range("H79").Value = 23799.8
$ws.Range("J79").Value = 23799.8
$ws.Range("L79").Value = 23799.8
$ws.Range("N79").Value = -25983.8

$ws.Range("H82").Value = 25981.908
$ws.Range("I82").Value = 4573.2856
$ws.Range("K82").Value = 4573.2856
$ws.Range("M82").Value = -4190.2856

$ws.Range("H85").Value = 25981.908
$ws.Range("I85").Value = 4573.2856
$ws.Range("K85").Value = 4573.2856
$ws.Range("M85").Value = -3247.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6728.115
$ws.Range("I31").Value = 1730.5714
$ws.Range("J31").Value = 8569.315000000001
$ws.Range("K31").Value = 1730.5714
$ws.Range("L31").Value = 8569.315000000001
$ws.Range("M31").Value = -1435.5714
$ws.Range("N31").Value = -9159.315000000001

$ws.Range("H34").Value = 6728.115
$ws.Range("I34").Value = 1730.5714
$ws.Range("J34").Value = 8569.315000000001
$ws.Range("K34").Value = 1730.5714
$ws.Range("L34").Value = 8569.315000000001
$ws.Range("M34").Value = -1528.5714
$ws.Range("N34").Value = -8973.315000000001

$ws.Range("H109").Value = 70259
$ws.Range("I109").Value = 70259
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 70259
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = $null
$ws.Range("M109").Value = -69219

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 845.9048
$ws.Range("J5").Value = 987.2
$ws.Range("L5").Value = 2961.6
$ws.Range("N5").Value = -3185.6

$ws.Range("H135").Value = 845.9048
$ws.Range("J135").Value = 987.2
$ws.Range("L135").Value = 8884.800000000001
$ws.Range("N135").Value = -13954.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 502.8
$ws.Range("I22").Value = 378.5
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 378.5
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = 150.5
$ws.Range("N22").Value = -2058

$ws.Range("H70").Value = 10805.4
$ws.Range("I70").Value = 10755.25
$ws.Range("J70").Value = 11006
$ws.Range("K70").Value = 10755.25
$ws.Range("L70").Value = 11006
$ws.Range("M70").Value = -10485.25
$ws.Range("N70").Value = -11546

$ws.Range("H73").Value = 10805.4
$ws.Range("I73").Value = 10755.25
$ws.Range("J73").Value = 11006
$ws.Range("K73").Value = 10755.25
$ws.Range("L73").Value = 11006
$ws.Range("M73").Value = -9819.25
$ws.Range("N73").Value = -12878

$ws.Range("H95").Value = 26622.25
$ws.Range("J95").Value = 26622.25
$ws.Range("L95").Value = 26622.25
$ws.Range("N95").Value = -32114.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8548.375
$ws.Range("I7").Value = 7608.5
$ws.Range("J7").Value = 9488.25
$ws.Range("K7").Value = 7608.5
$ws.Range("L7").Value = 9488.25
$ws.Range("M7").Value = -7496.5
$ws.Range("N7").Value = -9712.25

$ws.Range("H16").Value = 1049.2
$ws.Range("I16").Value = 1049.2
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1049.2
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -879.2
$ws.Range("N16").Value = $null

$ws.Range("H22").Value = 4000
$ws.Range("J22").Value = 4000
$ws.Range("L22").Value = 4000
$ws.Range("N22").Value = -4590

$ws.Range("H27").Value = 4000
$ws.Range("J27").Value = 4000
$ws.Range("L27").Value = 4000
$ws.Range("N27").Value = -4214

$ws.Range("H40").Value = 4500
$ws.Range("I40").Value = 4500
$ws.Range("K40").Value = 4500
$ws.Range("M40").Value = -4364

$ws.Range("H46").Value = 5297.8
$ws.Range("I46").Value = 3499.8
$ws.Range("K46").Value = 3499.8
$ws.Range("M46").Value = -3311.8

$ws.Range("H61").Value = 2435.2222
$ws.Range("I61").Value = 1322.2667
$ws.Range("K61").Value = 1322.2667
$ws.Range("M61").Value = -1120.2667

$ws.Range("H104").Value = 30000
$ws.Range("J104").Value = 30000
$ws.Range("L104").Value = 30000
$ws.Range("N104").Value = -36988

$ws.Range("H113").Value = 2435.2222
$ws.Range("I113").Value = 1322.2667
$ws.Range("K113").Value = 1322.2667
$ws.Range("M113").Value = 847.7333000000001

$ws.Range("H122").Value = 3247.25
$ws.Range("I122").Value = 3247.25
$ws.Range("K122").Value = 9741.75
$ws.Range("M122").Value = -7291.75

$ws.Range("H126").Value = 8548.375
$ws.Range("I126").Value = 7608.5
$ws.Range("J126").Value = 9488.25
$ws.Range("K126").Value = 22825.5
$ws.Range("L126").Value = 28464.75
$ws.Range("M126").Value = -20355.5
$ws.Range("N126").Value = -33404.75

$ws.Range("H136").Value = 8250
$ws.Range("I136").Value = 7500
$ws.Range("J136").Value = 9000
$ws.Range("K136").Value = 22500
$ws.Range("L136").Value = 27000
$ws.Range("M136").Value = -19950
$ws.Range("N136").Value = -32100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").Value = $null

$ws.Range("H132").Value = 2057.5715
$ws.Range("I132").Value = 2018.2667
$ws.Range("J132").Value = 2155.8333
$ws.Range("K132").Value = 6054.800099999999
$ws.Range("L132").Value = 6467.499899999999
$ws.Range("M132").Value = -3524.800099999999
$ws.Range("N132").Value = -11527.4999
